$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.321.53"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "3.376.64"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'570.75"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").Value = "'141.69"
$ws.Range("E6").Value = "  -4.63%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.378.20"
$ws.Range("E8").Value = "  -2.14%  "
$ws.Range("D9").Value = "'0.474"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'7.48"
$ws.Range("E10").Value = "  -4.05%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "'0.394"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "3.953.23"
$ws.Range("E13").Value = "  -2.29%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'28.10"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "'0.124"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "'0.0000171"
$ws.Range("E16").Value = "  -2.21%  "
$ws.Range("D17").Value = "3.380.56"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "60.534.42"
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("D19").Value = "'6.26"
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").Value = "'14.08"
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").Value = "'9.08"
$ws.Range("E21").Value = "  -3.87%  "
$ws.Range("D22").Value = "'388.88"
$ws.Range("D23").Value = "'0.560"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").Value = "'73.13"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  -4.74%  "
$ws.Range("D27").Value = "3.519.16"
$ws.Range("E27").Value = "  -2.15%  "
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "'7.40"
$ws.Range("E30").Value = "  -5.62%  "
$ws.Range("D31").Value = "'8.06"
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("D32").Value = "'2.17"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("E33").Value = "  -5.01%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'23.76"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").Value = "'6.95"
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("D37").Value = "3.405.98"
$ws.Range("E37").Value = "  -2.05%  "
$ws.Range("D38").Value = "'166.88"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").Value = "'4.99"
$ws.Range("E39").Value = "  -4.66%  "
$ws.Range("D40").Value = "'1.51"
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("D41").Value = "'0.0777"
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("D42").Value = "'26.86"
$ws.Range("E42").Value = "  +2.99%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.782"
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "'4.46"
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("D46").Value = "'41.35"
$ws.Range("E46").Value = "  -2.39%  "
$ws.Range("D47").Value = "'1.69"
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("D48").Value = "2.539.06"
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").Value = "'6.84"
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("D51").Value = "'22.95"
$ws.Range("E51").Value = "  -1.50%  "
